# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# de-de (and zh-cn) handback xliff files have been generated / handed back.

$wb = $excel.ActiveWorkbook

# Column widths end up quantized by the runtime to a 1/6-character grid
# (raw_width = round(ColumnWidth*6)/6 + 5/6). Pick the ColumnWidth values that
# land closest to the widths seen in the target workbook.
$wideColWidth  = 29.166666666666668   # -> raw width 30   (closest to 29.9777047293527)
$fortyColWidth = 39.166666666666664   # -> raw width 40   (exact)

# ---------------------------------------------------------------------------
# Overview sheet: the overall status moves from "Ready for handoff" to
# "Handed back: in sync with en-US" for both locales.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------------
# zh-cn sheet: record the already-generated handback artifacts.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Columns.Item(3).ColumnWidth  = $wideColWidth
$wsZhCn.Columns.Item(9).ColumnWidth  = $fortyColWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $fortyColWidth

$zhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b4df552e602c8248e829abb4cc8debd24516dacd/e2e/5c10fc32-de83-4758-a2f8-b35ef233937d.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhUrl, "", "", "5c10fc32-de83-4758-a2f8-b35ef233937d.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhUrl, "", "", "5c10fc32-de83-4758-a2f8-b35ef233937d.md") | Out-Null

$wsZhCn.Range("J2").Value = "5c10fc32-de83-4758-a2f8-b35ef233937d.546226525ee30590ae9a50a4ced6cfc33f38cde4.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "5c10fc32-de83-4758-a2f8-b35ef233937d.546226525ee30590ae9a50a4ced6cfc33f38cde4.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-08-25 13:06:06"
$wsZhCn.Range("K3").Value = "2016-08-25 13:06:06"

# ---------------------------------------------------------------------------
# de-de sheet: record the newly generated handback artifacts.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Columns.Item(3).ColumnWidth  = $wideColWidth
$wsDeDe.Columns.Item(9).ColumnWidth  = $fortyColWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $fortyColWidth

$deUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b4df552e602c8248e829abb4cc8debd24516dacd/e2e/5c10fc32-de83-4758-a2f8-b35ef233937d.md"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $deUrl, "", "", "5c10fc32-de83-4758-a2f8-b35ef233937d.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $deUrl, "", "", "5c10fc32-de83-4758-a2f8-b35ef233937d.md") | Out-Null

$wsDeDe.Range("J2").Value = "5c10fc32-de83-4758-a2f8-b35ef233937d.546226525ee30590ae9a50a4ced6cfc33f38cde4.de-de.xlf"
$wsDeDe.Range("J3").Value = "5c10fc32-de83-4758-a2f8-b35ef233937d.546226525ee30590ae9a50a4ced6cfc33f38cde4.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-25 13:06:17"
$wsDeDe.Range("K3").Value = "2016-08-25 13:06:17"
